# Update the "Förändrad" (Changed) date column (C) for existing data rows
# from 2024-10-15 (serial 45580) to 2024-10-16 (serial 45581), give row 29
# an explicit row height (matching the other data rows), and append a new
# data row (row 30) for case "A 45908-2024".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (rows 2 through 29) from 45580 -> 45581
for ($r = 2; $r -le 29; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45580) {
        $cell.Value = 45581
    }
}

# Row 29 now gets an explicit "15, custom height" row height, matching
# every other data row in the sheet.
$ws.Rows.Item(29).RowHeight = 15

# Append the new row 30.
$row = 30

$ws.Cells.Item($row, 1).Value = "A 45908-2024"

$ws.Cells.Item($row, 2).Value = 45580
$ws.Cells.Item($row, 2).NumberFormat = $ws.Cells.Item($row - 1, 2).NumberFormat

$ws.Cells.Item($row, 3).Value = 45581
$ws.Cells.Item($row, 3).NumberFormat = $ws.Cells.Item($row - 1, 3).NumberFormat

$ws.Cells.Item($row, 4).Value = "OKÄNT"
$ws.Cells.Item($row, 5).Value = "OKÄNT"

$ws.Cells.Item($row, 7).Value = 4.1
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 0
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0

# Column R carries the wrap-text style used by the rest of the table
# (style index 2 in the original file) even though the cell itself is blank.
$ws.Cells.Item($row, 18).WrapText = $true
